$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of workout data appended to the scoreboard log (rows 143-145)
$newRows = @(
    @{ A="Matt";   B=45474; C="Walk";    D=20; E=0.73;  F=66;  G=20; H=0;  I=0; J=0; K=0; L="Sauntering Hippo"; M=4 },
    @{ A="Steven"; B=45474; C="Workout"; D=23; E=0;     F=0;   G=23; H=0;  I=0; J=0; K=0; L="Mighty Monkey";    M=4 },
    @{ A="Matt";   B=45475; C="Ride";    D=57; E=14.77; F=863; G=9;  H=45; I=1; J=0; K=0; L="Sauntering Hippo"; M=4 }
)

$startRow = 143
$endRow = $startRow + $newRows.Count - 1

# Copy the date cell format (style index reused from the last existing
# date cell) down onto the new rows' date column before filling values.
$ws.Range("B142").Copy() | Out-Null
$ws.Range("B$startRow`:B$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
}

# Selection moves to A146, matching the post-edit cursor position
$ws.Range("A146").Select() | Out-Null
